$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): shift accuracy to B1, precision to C1, recall to D1
$ws.Range("B1").Value = "accuracy"
$ws.Range("C1").Value = "precision"
$ws.Range("D1").Value = "recall"

# Update row 2 (SVM_global): swap B2 and C2 values
$ws.Range("B2").Value = 92.85714285714286
$ws.Range("C2").Value = 86.22448979591836

# Update row 3 (AdaBoostClassifier_global) with new values
$ws.Range("B3").Value = 90.25974025974025
$ws.Range("C3").Value = 86.04761904761904
$ws.Range("D3").Value = 90.25974025974025
$ws.Range("E3").Value = 88.10336421257922
$ws.Range("F3").Value = 0.6223776223776223
